$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "노령연금 입금금액 = 330000"

$ws.Range("A3:G3").NumberFormat = "@"
$ws.Range("A3").Value = "325,920"
$ws.Range("B3").Value = "484,710"
$ws.Range("C3").Value = "643,500"
$ws.Range("D3").Value = "802,300"
$ws.Range("E3").Value = "961,090"
$ws.Range("F3").Value = "1,119,880"
$ws.Range("G3").Value = "1,278,680"
$ws.Range("A3:G3").ClearFormats()
